$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Split the run "e con possibilità di aggiungere textbox/buchi." into
#    three runs so that "textbox" is wrapped in spell-check proofErr
#    markers (spellStart/spellEnd), matching the author's edit. The
#    visible text is unchanged.
#
#    The whole owning paragraph is rewritten (rather than just the
#    affected run) so the untouched sibling runs keep their original
#    w:rsidR/w:rsidRPr attributes instead of being merged/stripped by
#    the host's run-normalisation pass.
# ----------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute(
    "Realizzare blocchi",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $pStart = $findRng.Start
    $pEnd = $findRng.End + 128
    $probe = $d.Range($pStart, $pEnd)
    # Walk forward from the match to the real end of the paragraph
    # (the character right after the trailing "." is the paragraph
    # mark itself).
    $paraEnd = $pStart
    $ch = ""
    do {
        $paraEnd = $paraEnd + 1
        $ch = $d.Range($paraEnd - 1, $paraEnd).Text
    } while ($ch -ne [char]13 -and $paraEnd -lt $d.Content.End)

    $wholePara1 = $d.Range($pStart, $paraEnd)

    $splitXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
<w:r w:rsidRPr="00C863FC"><w:rPr><w:strike/><w:lang w:val="it-IT"/></w:rPr><w:t>Realizzare blocchi</w:t></w:r>
<w:r w:rsidR="006522CD" w:rsidRPr="00C863FC"><w:rPr><w:strike/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> con estendibilit&#224; sulla base di quello che c&#8217;&#232; scritto dentro, dell&#8217;incastro </w:t></w:r>
<w:r w:rsidR="00FC44B0" w:rsidRPr="00C863FC"><w:rPr><w:strike/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">e con possibilit&#224; di aggiungere </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:strike/><w:lang w:val="it-IT"/></w:rPr><w:t>textbox</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:strike/><w:lang w:val="it-IT"/></w:rPr><w:t>/buchi.</w:t></w:r>
</w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@

    $wholePara1.InsertXML($splitXml)
}

# ----------------------------------------------------------------------
# 2) Insert a brand-new bullet paragraph right after the "Definire un
#    sistema di combinazione..." paragraph, with the "visitor" text
#    (proofErr spell markers around "visitor"), carrying over the
#    ListParagraph numbering/style and moving the _GoBack bookmark to
#    the end of the new paragraph (it currently sits at the end of the
#    "Definire..." paragraph).
# ----------------------------------------------------------------------
$defineRng = $d.Content
$foundDefine = $defineRng.Find.Execute(
    "Definire un sistema di combinazione dei blocchi e vedere dove questa strada ci porta: potrebbe portarci ad ottenere un parser implicito.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundDefine) {
    # The Find match covers only the sentence text itself; extending
    # the end by one character also captures the paragraph mark (and
    # the _GoBack bookmark sitting right before it), so the whole
    # paragraph (unchanged) plus the new paragraph can be written back
    # together in one shot.
    $pStart = $defineRng.Start
    $pEnd = $defineRng.End + 1
    $wholePara = $d.Range($pStart, $pEnd)

    $insertXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>Definire un sistema di combinazione dei blocchi e vedere dove questa strada ci porta: potrebbe portarci ad ottenere un parser implicito.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">Se volessimo estrarre un albero sintattico dal nostro sistema, si otterrebbe un albero non binario. In generale, non escluderei la possibilit&#224; di fare valutazioni con un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>visitor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>, anche se non occorre veramente avere pi&#249; strategie di valutazione. Bisogner&#224; infine prevedere un motore di valutazione unico che consenta la interpretazione degli script a divisione di tempo in modo equo.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@

    $wholePara.InsertXML($insertXml)
}
